$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report date range) ---
$ws.Range("A8").Value = "Volume 30   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/6/2023  Through  11/12/2023"

# --- Cells that change from numeric to text placeholders ("0" / "***.*") ---
# Use Copy() from stable source cells (row 14) that already hold the right
# shared-string + style combination, so style/type match exactly.
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("C14").Copy($ws.Range("D18"))
$ws.Range("E14").Copy($ws.Range("E18"))
$ws.Range("C14").Copy($ws.Range("F22"))
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("C14").Copy($ws.Range("D30"))
$ws.Range("E14").Copy($ws.Range("E30"))

# --- Plain numeric value updates ---
$ws.Range("M15").Value = 10.526315789473
$ws.Range("N15").Value = -34.375
$ws.Range("C16").Value = 5
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 23
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = -8
$ws.Range("I16").Value = 251
$ws.Range("J16").Value = 232
$ws.Range("K16").Value = 8.189655172413
$ws.Range("L16").Value = 32.105263157894
$ws.Range("M16").Value = -7.037037037037
$ws.Range("N16").Value = -58.305647840531
$ws.Range("C17").Value = 18
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 157.142857142857
$ws.Range("F17").Value = 38
$ws.Range("H17").Value = 22.580645161290
$ws.Range("I17").Value = 371
$ws.Range("J17").Value = 321
$ws.Range("K17").Value = 15.576323987538
$ws.Range("L17").Value = 47.808764940239
$ws.Range("M17").Value = 65.625
$ws.Range("N17").Value = 23.255813953488
$ws.Range("C18").Value = 2
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 120
$ws.Range("I18").Value = 209
$ws.Range("K18").Value = 84.955752212389
$ws.Range("L18").Value = 55.970149253731
$ws.Range("M18").Value = -37.611940298507
$ws.Range("N18").Value = -83.684621389539
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 110
$ws.Range("F19").Value = 65
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = 38.297872340425
$ws.Range("I19").Value = 595
$ws.Range("J19").Value = 532
$ws.Range("K19").Value = 11.842105263157
$ws.Range("L19").Value = 57.824933687002
$ws.Range("M19").Value = 52.956298200514
$ws.Range("N19").Value = 14.864864864864
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 25
$ws.Range("F20").Value = 37
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 452
$ws.Range("J20").Value = 308
$ws.Range("K20").Value = 46.753246753246
$ws.Range("L20").Value = 97.379912663755
$ws.Range("M20").Value = 128.282828282828
$ws.Range("N20").Value = -72.966507177033
$ws.Range("C21").Value = 56
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = 86.666666666666
$ws.Range("F21").Value = 174
$ws.Range("G21").Value = 147
$ws.Range("H21").Value = 18.367346938775
$ws.Range("I21").Value = 1904
$ws.Range("J21").Value = 1528
$ws.Range("K21").Value = 24.607329842931
$ws.Range("L21").Value = 57.87728026534
$ws.Range("M21").Value = 31.764705882352
$ws.Range("N21").Value = -56.903576278859
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = -45
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = 100
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 80
$ws.Range("I23").Value = 106
$ws.Range("J23").Value = 95
$ws.Range("K23").Value = 11.578947368421
$ws.Range("L23").Value = 34.177215189873
$ws.Range("M23").Value = 89.285714285714
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = -22.857142857142
$ws.Range("F24").Value = 134
$ws.Range("G24").Value = 126
$ws.Range("H24").Value = 6.349206349206
$ws.Range("I24").Value = 1410
$ws.Range("J24").Value = 1173
$ws.Range("K24").Value = 20.204603580562
$ws.Range("L24").Value = 60.227272727272
$ws.Range("M24").Value = 73.431734317343
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 220
$ws.Range("F25").Value = 43
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = 30.303030303030
$ws.Range("I25").Value = 464
$ws.Range("J25").Value = 460
$ws.Range("K25").Value = 0.869565217391
$ws.Range("L25").Value = 31.073446327683
$ws.Range("M25").Value = -17.142857142857
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -66.666666666666
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -75
$ws.Range("I27").Value = 61
$ws.Range("K27").Value = 29.787234042553
$ws.Range("L27").Value = 24.489795918367
$ws.Range("L28").Value = -27.777777777777
$ws.Range("N28").Value = -70.454545454545
$ws.Range("L29").Value = -14.285714285714
$ws.Range("N29").Value = -69.230769230769
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
